$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# List1 : drop the "contextual spacing" paragraph flag.
# (Style.NoSpaceBetweenParagraphsOfSameStyle is the COM-level name for
#  <w:contextualSpacing/>; the style is recreated so the flag element is
#  omitted entirely rather than written out as an explicit "false".)
# ---------------------------------------------------------------------------
$d.Styles("List1").Delete()
$list1 = $d.Styles.Add("List1", 1)
$list1.NameLocal = "List 1"
$list1.BaseStyle = $d.Styles("List")
$list1.LinkStyle = $d.Styles("List1Char")
$list1.ParagraphFormat.SpaceAfter = 8
$list1.ParagraphFormat.LineSpacingRule = 0
$list1.ParagraphFormat.LineSpacing = 12
$list1.Font.Name = "Times New Roman"
$list1.Font.NameBi = "Times New Roman"
$list1.Font.Size = 12
$list1.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# List6 : no longer based on Heading4 - now a standalone style with its own
# run properties (Times New Roman 12pt, not bold) and no outline level /
# contextual spacing.
# ---------------------------------------------------------------------------
$d.Styles("List6").Delete()
$list6 = $d.Styles.Add("List6", 1)
$list6.NameLocal = "List 6"
$list6.LinkStyle = $d.Styles("List6Char")
$list6.ParagraphFormat.LineSpacingRule = 5
$list6.ParagraphFormat.LineSpacing = 12.95
$list6.ParagraphFormat.FirstLineIndent = 108
$list6.Font.Name = "Times New Roman"
$list6.Font.NameBi = "Times New Roman"
$list6.Font.Size = 12
$list6.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# List7 : no longer based on Heading4; keeps contextual spacing, becomes
# bold, gains its own run properties.
# ---------------------------------------------------------------------------
$d.Styles("List7").Delete()
$list7 = $d.Styles.Add("List7", 1)
$list7.NameLocal = "List 7"
$list7.LinkStyle = $d.Styles("List7Char")
$list7.ParagraphFormat.SpaceAfter = 8
$list7.ParagraphFormat.LineSpacingRule = 5
$list7.ParagraphFormat.LineSpacing = 12.95
$list7.ParagraphFormat.FirstLineIndent = 162
$list7.NoSpaceBetweenParagraphsOfSameStyle = $true
$list7.Font.Name = "Times New Roman"
$list7.Font.NameBi = "Times New Roman"
$list7.Font.Size = 12
$list7.Font.SizeBi = 12
$list7.Font.Bold = $true

# List7Char : the linked character style keeps its Heading4Char base, but
# flips from "not bold" to "bold".
$d.Styles("List7Char").Font.Bold = $true

# ---------------------------------------------------------------------------
# List8 : no longer based on Heading4; gains its own run properties, no
# outline level / contextual spacing.
# ---------------------------------------------------------------------------
$d.Styles("List8").Delete()
$list8 = $d.Styles.Add("List8", 1)
$list8.NameLocal = "List 8"
$list8.LinkStyle = $d.Styles("List8Char")
$list8.ParagraphFormat.SpaceAfter = 8
$list8.ParagraphFormat.LineSpacingRule = 5
$list8.ParagraphFormat.LineSpacing = 12.95
$list8.ParagraphFormat.FirstLineIndent = 180
$list8.Font.Name = "Times New Roman"
$list8.Font.NameBi = "Times New Roman"
$list8.Font.Size = 12
$list8.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# List3change / List4change : drop the "contextual spacing" paragraph flag.
# ---------------------------------------------------------------------------
$d.Styles("List3change").Delete()
$list3change = $d.Styles.Add("List3change", 1)
$list3change.NameLocal = "List 3_change"
$list3change.BaseStyle = $d.Styles("List1")
$list3change.LinkStyle = $d.Styles("List3changeChar")
$list3change.ParagraphFormat.SpaceAfter = 0
$list3change.ParagraphFormat.LeftIndent = 54

$d.Styles("List4change").Delete()
$list4change = $d.Styles.Add("List4change", 1)
$list4change.NameLocal = "List 4_change"
$list4change.BaseStyle = $d.Styles("List1")
$list4change.LinkStyle = $d.Styles("List4changeChar")
$list4change.ParagraphFormat.SpaceAfter = 0
$list4change.ParagraphFormat.LeftIndent = 72

Write-Output "styles updated"
